$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 97.57717366666668
$ws.Range("H2").Value = 292.731521
$ws.Range("I2").Value = 0.3532166605548384
$ws.Range("J2").Value = 0.3532166605548384
$ws.Range("M2").Value = 21.106626
$ws.Range("N2").Value = 63.319878
$ws.Range("O2").Value = 0.3106213714361249
$ws.Range("P2").Value = 0.3106213714361249
$ws.Range("Q2").Value = 2059.524910719383
$ws.Range("R2").Value = 18535.72419647444
$ws.Range("S2").Value = 0.1097166435156321
$ws.Range("T2").Value = 0.1097166435156321
$ws.Range("G3").Value = 97.57717366666668
$ws.Range("H3").Value = 292.731521
$ws.Range("I3").Value = 0.3532166605548384
$ws.Range("J3").Value = 0.3532166605548384
$ws.Range("O3").Value = 0.03085709917216154
$ws.Range("P3").Value = 0.03085709917216154
$ws.Range("Q3").Value = 204.5930198678346
$ws.Range("R3").Value = 1841.337178810511
$ws.Range("S3").Value = 0.01089924152400037
$ws.Range("T3").Value = 0.01089924152400037
$ws.Range("G4").Value = 97.57717366666668
$ws.Range("H4").Value = 292.731521
$ws.Range("I4").Value = 0.3532166605548384
$ws.Range("J4").Value = 0.3532166605548384
$ws.Range("O4").Value = 0.6585215293917135
$ws.Range("P4").Value = 0.6585215293917135
$ws.Range("Q4").Value = 4366.220803664672
$ws.Range("R4").Value = 39295.98723298205
$ws.Range("S4").Value = 0.2326007755152059
$ws.Range("T4").Value = 0.2326007755152059
$ws.Range("I5").Value = 0.5533024543641269
$ws.Range("J5").Value = 0.5533024543641269
$ws.Range("M5").Value = 21.106626
$ws.Range("N5").Value = 63.319878
$ws.Range("O5").Value = 0.3106213714361249
$ws.Range("P5").Value = 0.3106213714361249
$ws.Range("Q5").Value = 3226.179043013106
$ws.Range("R5").Value = 29035.61138711796
$ws.Range("S5").Value = 0.1718675671935591
$ws.Range("T5").Value = 0.171867567193559
$ws.Range("I6").Value = 0.5533024543641269
$ws.Range("J6").Value = 0.5533024543641269
$ws.Range("O6").Value = 0.03085709917216154
$ws.Range("P6").Value = 0.03085709917216154
$ws.Range("S6").Value = 0.01707330870651425
$ws.Range("T6").Value = 0.01707330870651425
$ws.Range("I7").Value = 0.5533024543641269
$ws.Range("J7").Value = 0.5533024543641269
$ws.Range("O7").Value = 0.6585215293917135
$ws.Range("P7").Value = 0.6585215293917135
$ws.Range("S7").Value = 0.3643615784640536
$ws.Range("T7").Value = 0.3643615784640536
$ws.Range("I8").Value = 0.09348088508103472
$ws.Range("J8").Value = 0.09348088508103473
$ws.Range("M8").Value = 21.106626
$ws.Range("N8").Value = 63.319878
$ws.Range("O8").Value = 0.3106213714361249
$ws.Range("P8").Value = 0.3106213714361249
$ws.Range("Q8").Value = 545.0654881286281
$ws.Range("R8").Value = 4905.589393157652
$ws.Range("S8").Value = 0.0290371607269338
$ws.Range("T8").Value = 0.0290371607269338
$ws.Range("I9").Value = 0.09348088508103472
$ws.Range("J9").Value = 0.09348088508103473
$ws.Range("O9").Value = 0.03085709917216154
$ws.Range("P9").Value = 0.03085709917216154
$ws.Range("S9").Value = 0.002884548941646925
$ws.Range("T9").Value = 0.002884548941646925
$ws.Range("I10").Value = 0.09348088508103472
$ws.Range("J10").Value = 0.09348088508103473
$ws.Range("O10").Value = 0.6585215293917135
$ws.Range("P10").Value = 0.6585215293917135
$ws.Range("S10").Value = 0.061559175412454
$ws.Range("T10").Value = 0.06155917541245401
